$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each target cell is plain text (inlineStr) in the canonical file, even when
# the text looks like a number (e.g. "318.38") or like a multi-dot thousands
# value (e.g. "44.067.40"). Excel auto-converts numeric-looking strings typed
# into a cell to real numbers, so force the cell to Text format first, assign
# the literal string, then restore the "Normal" style so no visible formatting
# change is left behind.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "44.067.40"
Set-TextValue "E2" "  +2.08%  "
Set-TextValue "D3" "2.256.09"
Set-TextValue "E3" "  +1.47%  "
Set-TextValue "E4" "  +0.09%  "
Set-TextValue "D5" "318.38"
Set-TextValue "E5" "  -0.35%  "
Set-TextValue "D6" "100.91"
Set-TextValue "E6" "  +2.21%  "
Set-TextValue "D7" "0.577"
Set-TextValue "E7" "  -0.65%  "
Set-TextValue "E8" "  +0.16%  "
Set-TextValue "D9" "0.551"
Set-TextValue "E9" "  -2.45%  "
Set-TextValue "D10" "37.19"
Set-TextValue "E10" "  +1.22%  "
Set-TextValue "E11" "  +1.11%  "
Set-TextValue "D12" "7.58"
Set-TextValue "E12" "  -0.66%  "
Set-TextValue "E13" "  -1.28%  "
Set-TextValue "D14" "2.601.11"
Set-TextValue "E14" "  +1.62%  "
Set-TextValue "D15" "14.53"
Set-TextValue "E15" "  +1.68%  "
Set-TextValue "E16" "  -0.30%  "
Set-TextValue "D17" "2.258.09"
Set-TextValue "E17" "  +1.72%  "
Set-TextValue "D18" "43.953.30"
Set-TextValue "E18" "  +2.15%  "
Set-TextValue "D19" "13.40"
Set-TextValue "E19" "  -2.62%  "
Set-TextValue "E20" "  +2.36%  "
Set-TextValue "D21" "6.45"
Set-TextValue "E21" "  -1.27%  "
Set-TextValue "D22" "65.61"
Set-TextValue "E22" "  +0.69%  "
Set-TextValue "D23" "3.10"
Set-TextValue "E23" "  -4.26%  "
Set-TextValue "D24" "234.78"
Set-TextValue "E24" "  -0.49%  "
Set-TextValue "E25" "  -5.67%  "
Set-TextValue "E26" "  +0.33%  "
Set-TextValue "D27" "10.76"
Set-TextValue "E27" "  +7.45%  "
Set-TextValue "D28" "38.74"
Set-TextValue "E28" "  +6.36%  "
Set-TextValue "E29" "  -0.85%  "
Set-TextValue "D30" "6.17"
Set-TextValue "E30" "  -3.08%  "
Set-TextValue "D31" "161.04"
Set-TextValue "E31" "  +3.48%  "
Set-TextValue "D32" "20.18"
Set-TextValue "E32" "  -0.27%  "
Set-TextValue "E33" "  -1.77%  "
Set-TextValue "E34" "  +1.01%  "
Set-TextValue "D35" "1.98"
Set-TextValue "E35" "  +7.42%  "
Set-TextValue "B36" "LidoDAOToken"
Set-TextValue "C36" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D36" "3.10"
Set-TextValue "E36" "  -6.95%  "
Set-TextValue "B37" "Kaspa"
Set-TextValue "C37" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D37" "0.113"
Set-TextValue "E37" "  +7.17%  "
Set-TextValue "E38" "  -1.45%  "
Set-TextValue "D39" "16.51"
Set-TextValue "E39" "  +18.02%  "
Set-TextValue "B40" "NEARProtocol"
Set-TextValue "C40" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D40" "3.67"
Set-TextValue "E40" "  -0.07%  "
Set-TextValue "B41" "RenderToken"
Set-TextValue "C41" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D41" "4.16"
Set-TextValue "E41" "  -5.00%  "
Set-TextValue "D42" "0.0317"
Set-TextValue "E42" "  -0.64%  "
Set-TextValue "E43" "  +0.10%  "
Set-TextValue "D44" "1.776.54"
Set-TextValue "E44" "  +2.63%  "
Set-TextValue "D45" "0.197"
Set-TextValue "E45" "  -2.59%  "
Set-TextValue "D46" "74.79"
Set-TextValue "E46" "  +1.18%  "
Set-TextValue "D47" "5.19"
Set-TextValue "E47" "  -1.32%  "
Set-TextValue "D48" "81.36"
Set-TextValue "E48" "  -3.28%  "
Set-TextValue "E49" "  +1.42%  "
Set-TextValue "D50" "57.97"
Set-TextValue "E50" "  +0.49%  "
Set-TextValue "D51" "1.67"
Set-TextValue "E51" "  +5.31%  "
